$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shear center reference point added: update row 2 (A2, B2) values
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 0.5
